$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Coin/Link/Price/Volume columns keep their original text formatting
# (NumberFormat "@" stops Excel from auto-coercing numeric-looking / percent
# strings into actual numbers when the value is assigned).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.15%"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.92"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.47%"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.154"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.20%"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08147"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.27%"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.947"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.56%"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.156"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.02%"

# Row 8
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9299"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.17%"

# Row 9
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1437"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.46%"

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1926"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.46%"

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09142"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.23%"

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03516"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.48%"

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09786"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.81%"

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001392"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.63%"

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005855"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-6.08%"

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.918"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.08%"

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.243"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.40%"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.325"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.31%"

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.57%"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1313"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.69%"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.632"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-4.01%"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2427"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.13%"

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.04%"

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.92%"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004372"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.75%"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.15%"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004005"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-9.96%"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02051"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.19%"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05069"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.14%"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007402"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.54%"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009864"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.61%"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1365"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.53%"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002130"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.16%"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009379"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.24%"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006370"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.92%"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.03%"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002717"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-18.84%"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.03%"

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.03%"
